$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.22224168450624
$ws.Range("C2").Value = 10.71512280286241
$ws.Range("D2").Value = 6.366890670965001
$ws.Range("E2").Value = 13.22350241928196
$ws.Range("F2").Value = 30.5146448945708
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("K2").Value = 8.61094144185912
$ws.Range("L2").Value = 9.96923161932693
$ws.Range("M2").Value = 14.07003342166148
$ws.Range("O2").Value = 27.3772848623333

$ws.Range("B3").Value = 12.01219342924657
$ws.Range("C3").Value = 10.71862745352171
$ws.Range("D3").Value = 6.331521204465345
$ws.Range("E3").Value = 13.2538932527781
$ws.Range("F3").Value = 30.55624990641657
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("K3").Value = 8.448034244293421
$ws.Range("L3").Value = 9.977100882347445
$ws.Range("M3").Value = 14.04258394266601
$ws.Range("O3").Value = 27.45032074583864

$ws.Range("B4").Value = 11.88372986712436
$ws.Range("C4").Value = 10.72113088122678
$ws.Range("D4").Value = 6.309310850992487
$ws.Range("E4").Value = 13.27431673644077
$ws.Range("F4").Value = 30.58872595969658
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("K4").Value = 8.347819480564281
$ws.Range("L4").Value = 9.983261192141244
$ws.Range("M4").Value = 14.02779152086777
$ws.Range("O4").Value = 27.50030722770337

$ws.Range("B5").Value = 11.83157727433648
$ws.Range("C5").Value = 10.72223984542566
$ws.Range("D5").Value = 6.300138251206354
$ws.Range("E5").Value = 13.28308310912181
$ws.Range("F5").Value = 30.6037005553944
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("K5").Value = 8.306988566370848
$ws.Range("L5").Value = 9.986106150725195
$ws.Range("M5").Value = 14.02228649506488
$ws.Range("O5").Value = 27.52196834428897

$ws.Range("B6").Value = 11.82293134906925
$ws.Range("C6").Value = 10.72242936282232
$ws.Range("D6").Value = 6.298607866714666
$ws.Range("E6").Value = 13.28456555939736
$ws.Range("F6").Value = 30.60629210701302
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("K6").Value = 8.300210699100047
$ws.Range("L6").Value = 9.986598775497646
$ws.Range("M6").Value = 14.02140410497081
$ws.Range("O6").Value = 27.52564308733063

$ws.Range("B7").Value = 11.88302562718395
$ws.Range("C7").Value = 10.72114547704487
$ws.Range("D7").Value = 6.309187635750489
$ws.Range("E7").Value = 13.27443316625387
$ws.Range("F7").Value = 30.58892086934948
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("K7").Value = 8.347268717140507
$ws.Range("L7").Value = 9.983298204909797
$ws.Range("M7").Value = 14.02771515474147
$ws.Range("O7").Value = 27.50059413150322

$ws.Range("B8").Value = 12.14975354108283
$ws.Range("C8").Value = 10.71625847630712
$ws.Range("D8").Value = 6.354799082283151
$ws.Range("E8").Value = 13.23361532287097
$ws.Range("F8").Value = 30.52755064520935
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("K8").Value = 8.554842898146557
$ws.Range("L8").Value = 9.971669521759262
$ws.Range("M8").Value = 14.06014379282183
$ws.Range("O8").Value = 27.40139932169658

$ws.Range("B9").Value = 12.67380393761567
$ws.Range("C9").Value = 10.70944754489793
$ws.Range("D9").Value = 6.44022728213915
$ws.Range("E9").Value = 13.16755344902226
$ws.Range("F9").Value = 30.46227634945114
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("K9").Value = 8.958046491994546
$ws.Range("L9").Value = 9.959383263496676
$ws.Range("M9").Value = 14.13987959463981
$ws.Range("O9").Value = 27.2477569736742

$ws.Range("B10").Value = 13.05548486233775
$ws.Range("C10").Value = 10.70611140046324
$ws.Range("D10").Value = 6.500411454563333
$ws.Range("E10").Value = 13.12752766114887
$ws.Range("F10").Value = 30.44797197569103
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("K10").Value = 9.248909291411383
$ws.Range("L10").Value = 9.956733812941463
$ws.Range("M10").Value = 14.2079949380808
$ws.Range("O10").Value = 27.15989728792775

$ws.Range("B11").Value = 13.2276011006314
$ws.Range("C11").Value = 10.70495117608616
$ws.Range("D11").Value = 6.527203043353504
$ws.Range("E11").Value = 13.11116354257835
$ws.Range("F11").Value = 30.44877527234219
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("K11").Value = 9.379465579040707
$ws.Range("L11").Value = 9.956904597120303
$ws.Range("M11").Value = 14.24097991776738
$ws.Range("O11").Value = 27.12537982569004

$ws.Range("B12").Value = 13.29249677107259
$ws.Range("C12").Value = 10.7045628310169
$ws.Range("D12").Value = 6.537261627118469
$ws.Range("E12").Value = 13.1052317390207
$ws.Range("F12").Value = 30.45012969705661
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("K12").Value = 9.428604362175465
$ws.Range("L12").Value = 9.957166268275017
$ws.Range("M12").Value = 14.25375122667105
$ws.Range("O12").Value = 27.11309404327483

$ws.Range("B13").Value = 13.27853380497551
$ws.Range("C13").Value = 10.70464420569445
$ws.Range("D13").Value = 6.535099231063334
$ws.Range("E13").Value = 13.10649747955498
$ws.Range("F13").Value = 30.44979131165857
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("K13").Value = 9.418035513045261
$ws.Range("L13").Value = 9.957101165821838
$ws.Range("M13").Value = 14.25098832859558
$ws.Range("O13").Value = 27.11570506113389

$ws.Range("B14").Value = 13.23294609163703
$ws.Range("C14").Value = 10.70491820666787
$ws.Range("D14").Value = 6.528032322943575
$ws.Range("E14").Value = 13.11067022111966
$ws.Range("F14").Value = 30.44886566173042
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("K14").Value = 9.383514530298305
$ws.Range("L14").Value = 9.956922182075733
$ws.Range("M14").Value = 14.24202503739328
$ws.Range("O14").Value = 27.12435332187504

$ws.Range("B15").Value = 13.20498387026178
$ws.Range("C15").Value = 10.70509267127782
$ws.Range("D15").Value = 6.52369225962649
$ws.Range("E15").Value = 13.11326064149964
$ws.Range("F15").Value = 30.44843540383717
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("K15").Value = 9.362329002533761
$ws.Range("L15").Value = 9.956838177261638
$ws.Range("M15").Value = 14.23657109862348
$ws.Range("O15").Value = 27.12975293272837

$ws.Range("B16").Value = 13.04420072950741
$ws.Range("C16").Value = 10.70619438257919
$ws.Range("D16").Value = 6.498648547486701
$ws.Range("E16").Value = 13.12863417991722
$ws.Range("F16").Value = 30.44806651424304
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("K16").Value = 9.240337650294972
$ws.Range("L16").Value = 9.956750268161379
$ws.Range("M16").Value = 14.20587893920755
$ws.Range("O16").Value = 27.16226289018374

$ws.Range("B17").Value = 12.94513198690425
$ws.Range("C17").Value = 10.70696151038304
$ws.Range("D17").Value = 6.483133058291622
$ws.Range("E17").Value = 13.1385374526119
$ws.Range("F17").Value = 30.44971210955058
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("K17").Value = 9.165015078039893
$ws.Range("L17").Value = 9.957048176170449
$ws.Range("M17").Value = 14.18755751333152
$ws.Range("O17").Value = 27.18360376307114

$ws.Range("B18").Value = 12.88801153824389
$ws.Range("C18").Value = 10.7074364056857
$ws.Range("D18").Value = 6.474153901935469
$ws.Range("E18").Value = 13.14440711499496
$ws.Range("F18").Value = 30.45134680995907
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("K18").Value = 9.121528648583569
$ws.Range("L18").Value = 9.957349100582205
$ws.Range("M18").Value = 14.17720809553754
$ws.Range("O18").Value = 27.19639140283586

$ws.Range("B19").Value = 12.86864969392406
$ws.Range("C19").Value = 10.7076029906799
$ws.Range("D19").Value = 6.47110433352762
$ws.Range("E19").Value = 13.14642429717951
$ws.Range("F19").Value = 30.45201850852506
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("K19").Value = 9.10677839674606
$ws.Range("L19").Value = 9.957473268534532
$ws.Range("M19").Value = 14.1737365527549
$ws.Range("O19").Value = 27.2008091272723

$ws.Range("B20").Value = 12.95569285131767
$ws.Range("C20").Value = 10.70687636706612
$ws.Range("D20").Value = 6.484790424824008
$ws.Range("E20").Value = 13.13746527075365
$ws.Range("F20").Value = 30.44946571215555
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("K20").Value = 9.173050501784568
$ws.Range("L20").Value = 9.957003058855268
$ws.Range("M20").Value = 14.18948839141428
$ws.Range("O20").Value = 27.18127889273386

$ws.Range("B21").Value = 13.24634442613805
$ws.Range("C21").Value = 10.70483634471428
$ws.Range("D21").Value = 6.530110420706
$ws.Range("E21").Value = 13.10943739729748
$ws.Range("F21").Value = 30.44910905606689
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("K21").Value = 9.393662685481502
$ws.Range("L21").Value = 9.956969414772205
$ws.Range("M21").Value = 14.24465021355486
$ws.Range("O21").Value = 27.12179179611969

$ws.Range("B22").Value = 13.43463677068183
$ws.Range("C22").Value = 10.70380019956016
$ws.Range("D22").Value = 6.55922254014723
$ws.Range("E22").Value = 13.09266368524717
$ws.Range("F22").Value = 30.45499654809129
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("K22").Value = 9.536076623143638
$ws.Range("L22").Value = 9.958095229306565
$ws.Range("M22").Value = 14.28233389400155
$ws.Range("O22").Value = 27.08749104046231

$ws.Range("B23").Value = 13.33431417609994
$ws.Range("C23").Value = 10.7043261502395
$ws.Range("D23").Value = 6.543732064599276
$ws.Range("E23").Value = 13.10147491635816
$ws.Range("F23").Value = 30.4512947666089
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("K23").Value = 9.460244357410854
$ws.Range("L23").Value = 9.957389644434372
$ws.Range("M23").Value = 14.26207439208155
$ws.Range("O23").Value = 27.10537871113882

$ws.Range("B24").Value = 12.95091879205755
$ws.Range("C24").Value = 10.70691475485048
$ws.Range("D24").Value = 6.484041313288164
$ws.Range("E24").Value = 13.13794945550392
$ws.Range("F24").Value = 30.44957496331171
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("K24").Value = 9.169418251477552
$ws.Range("L24").Value = 9.957023052474675
$ws.Range("M24").Value = 14.18861486835741
$ws.Range("O24").Value = 27.1823283516471

$ws.Range("B25").Value = 12.53234256692108
$ws.Range("C25").Value = 10.71099558358733
$ws.Range("D25").Value = 6.417560273952438
$ws.Range("E25").Value = 13.18392945353029
$ws.Range("F25").Value = 30.47402808996053
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("K25").Value = 8.849701370354127
$ws.Range("L25").Value = 9.961584333315534
$ws.Range("M25").Value = 14.11661110055496
$ws.Range("O25").Value = 27.28493404926462
